# Scheduled market-data refresh: update computed price/profit columns
# (currentAveragePrice, currentAveragePriceNQ/HQ, LevePriceNQ/HQ,
# LeveProfitNQ/HQ) on each class sheet to the latest scrape. A few rows'
# HQ-profit columns have no HQ recipe/listing this cycle, so those cells
# are cleared outright rather than zeroed.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 2800
$ws.Range("I51").Value = 2700
$ws.Range("K51").Value = 2700
$ws.Range("M51").Value = -2216
$ws.Range("H113").Value = 3114
$ws.Range("I113").Value = 2336.6667
$ws.Range("J113").Value = 3502.6667
$ws.Range("K113").Value = 2336.6667
$ws.Range("L113").Value = 3502.6667
$ws.Range("M113").Value = 917.3332999999998
$ws.Range("N113").Value = -10010.6667
$ws.Range("H129").Value = 1103.3766
$ws.Range("J129").Value = 1182.7142
$ws.Range("L129").Value = 3548.1426
$ws.Range("N129").Value = -13548.1426
$ws.Range("H138").Value = 3310.7812
$ws.Range("J138").Value = 4957.6772
$ws.Range("L138").Value = 14873.0316
$ws.Range("N138").Value = -25153.0316
$ws.Range("H141").Value = 5565.615
$ws.Range("I141").Value = 1813.7391
$ws.Range("J141").Value = 34330
$ws.Range("K141").Value = 5441.2173
$ws.Range("L141").Value = 102990
$ws.Range("M141").Value = -261.2173000000003
$ws.Range("N141").Value = -113350

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("M4").ClearContents()
$ws.Range("H5").Value = 2643.625
$ws.Range("I5").Value = 2643.625
$ws.Range("K5").Value = 2643.625
$ws.Range("M5").Value = -2531.625
$ws.Range("H74").Value = 1489.4706
$ws.Range("I74").Value = 1355.5714
$ws.Range("J74").Value = 1782.375
$ws.Range("K74").Value = 1355.5714
$ws.Range("L74").Value = 1782.375
$ws.Range("M74").Value = -481.5714
$ws.Range("N74").Value = -3530.375
$ws.Range("H77").Value = 1489.4706
$ws.Range("I77").Value = 1355.5714
$ws.Range("J77").Value = 1782.375
$ws.Range("K77").Value = 6777.857
$ws.Range("L77").Value = 8911.875
$ws.Range("M77").Value = -2409.857
$ws.Range("N77").Value = -17647.875

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 2643.625
$ws.Range("I4").Value = 2643.625
$ws.Range("K4").Value = 2643.625
$ws.Range("M4").Value = -2528.625
$ws.Range("H14").Value = 46500
$ws.Range("J14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("N14").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 90
$ws.Range("J7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("N7").ClearContents()
$ws.Range("H97").Value = 29500
$ws.Range("J97").Value = 29500
$ws.Range("L97").Value = 29500
$ws.Range("N97").Value = -31482
$ws.Range("H105").Value = 12777.777
$ws.Range("I105").Value = 15285.714
$ws.Range("J105").Value = 4000
$ws.Range("K105").Value = 15285.714
$ws.Range("L105").Value = 4000
$ws.Range("M105").Value = -13538.714
$ws.Range("N105").Value = -7494
$ws.Range("H132").Value = 3521.3076
$ws.Range("I132").Value = 2907.8
$ws.Range("K132").Value = 8723.400000000001
$ws.Range("M132").Value = -6193.400000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 369.7
$ws.Range("I5").Value = 388.55554
$ws.Range("J5").Value = 200
$ws.Range("K5").Value = 1165.66662
$ws.Range("L5").Value = 600
$ws.Range("M5").Value = -1053.66662
$ws.Range("N5").Value = -824
$ws.Range("H113").Value = 244604.12
$ws.Range("I113").Value = 588793.1
$ws.Range("J113").Value = 803.5833
$ws.Range("K113").Value = 1766379.3
$ws.Range("L113").Value = 2410.7499
$ws.Range("M113").Value = -1764209.3
$ws.Range("N113").Value = -6750.7499
$ws.Range("H114").Value = 1436.9412
$ws.Range("I114").Value = 409.33334
$ws.Range("J114").Value = 1997.4546
$ws.Range("K114").Value = 1228.00002
$ws.Range("L114").Value = 5992.3638
$ws.Range("M114").Value = 2025.99998
$ws.Range("N114").Value = -12500.3638
$ws.Range("H120").Value = 9106
$ws.Range("I120").Value = 6382.5
$ws.Range("K120").Value = 19147.5
$ws.Range("M120").Value = -14309.5
$ws.Range("H131").Value = 2833.255
$ws.Range("I131").Value = 484.6111
$ws.Range("J131").Value = 4114.3335
$ws.Range("K131").Value = 1453.8333
$ws.Range("L131").Value = 12343.0005
$ws.Range("M131").Value = 3586.1667
$ws.Range("N131").Value = -22423.0005
$ws.Range("H132").Value = 1270.4783
$ws.Range("I132").Value = 904.7857
$ws.Range("J132").Value = 1839.3334
$ws.Range("K132").Value = 8143.071300000001
$ws.Range("L132").Value = 16554.0006
$ws.Range("M132").Value = -5613.071300000001
$ws.Range("N132").Value = -21614.0006
$ws.Range("H135").Value = 369.7
$ws.Range("I135").Value = 388.55554
$ws.Range("J135").Value = 200
$ws.Range("K135").Value = 3496.99986
$ws.Range("L135").Value = 1800
$ws.Range("M135").Value = -961.9998599999999
$ws.Range("N135").Value = -6870

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 41.875
$ws.Range("I2").Value = 41.2
$ws.Range("K2").Value = 41.2
$ws.Range("M2").Value = 71.8
$ws.Range("H33").Value = 0
$ws.Range("I33").Value = 0
$ws.Range("K33").Value = 0
$ws.Range("M33").ClearContents()
$ws.Range("H51").Value = 35040
$ws.Range("J51").Value = 35040
$ws.Range("L51").Value = 35040
$ws.Range("N51").Value = -36058
$ws.Range("H122").Value = 3389
$ws.Range("I122").Value = 2778
$ws.Range("J122").Value = 4000
$ws.Range("K122").Value = 8334
$ws.Range("L122").Value = 12000
$ws.Range("M122").Value = -5884
$ws.Range("N122").Value = -16900
$ws.Range("H123").Value = 22254.38
$ws.Range("J123").Value = 22254.38
$ws.Range("L123").Value = 22254.38
$ws.Range("N123").Value = -27154.38

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 6084
$ws.Range("I9").Value = 443.33334
$ws.Range("J9").Value = 11724.667
$ws.Range("K9").Value = 443.33334
$ws.Range("L9").Value = 11724.667
$ws.Range("M9").Value = -219.33334
$ws.Range("N9").Value = -12172.667
$ws.Range("H55").Value = 275
$ws.Range("I55").Value = 268.875
$ws.Range("J55").Value = 287.25
$ws.Range("K55").Value = 268.875
$ws.Range("L55").Value = 287.25
$ws.Range("M55").Value = -95.875
$ws.Range("N55").Value = -633.25
$ws.Range("H61").Value = 4246.6665
$ws.Range("I61").Value = 2870
$ws.Range("K61").Value = 2870
$ws.Range("M61").Value = -2668
$ws.Range("H113").Value = 4246.6665
$ws.Range("I113").Value = 2870
$ws.Range("K113").Value = 2870
$ws.Range("M113").Value = -700
$ws.Range("H122").Value = 12505176
$ws.Range("J122").Value = 18187136
$ws.Range("L122").Value = 54561408
$ws.Range("N122").Value = -54566308

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 5740
$ws.Range("I62").Value = 7000
$ws.Range("J62").Value = 4900
$ws.Range("K62").Value = 7000
$ws.Range("L62").Value = 4900
$ws.Range("M62").Value = -6376
$ws.Range("N62").Value = -6148
$ws.Range("H65").Value = 5740
$ws.Range("I65").Value = 7000
$ws.Range("J65").Value = 4900
$ws.Range("K65").Value = 35000
$ws.Range("L65").Value = 24500
$ws.Range("M65").Value = -31880
$ws.Range("N65").Value = -30740
$ws.Range("H123").Value = 41952.81
$ws.Range("J123").Value = 41952.81
$ws.Range("L123").Value = 41952.81
$ws.Range("N123").Value = -51752.81
